$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "BarClose(timeframe) - 1"
$ws.Range("A12").Value = "NeuerDreierNeu(ID) - 2"
$ws.Range("A15").Value = "DreierKaputt(ID) - 3"

$ws.Range("E6").Select()
